# vocabulary.xlsx update: new .ttl generated from the Google sheet added a
# "covid19" PREFIX declaration plus a new "spatial scope" (covid19:10001)
# concept that "diary cattle" (id-amr:10001) is now a narrower term of.
#
# Net effect vs. the original sheet:
#   * a new row is inserted at row 8  -> PREFIX / covid19 / http://purl.org/zonmw/covid19/
#   * a new row is inserted at row 19 -> covid19:10001 / spatial scope
#   * the (now shifted) "diary cattle" row (row 20) gets a skos:broader
#     reference to covid19:10001 in column E, and a single space in column G
#   * everything else simply shifts down by two rows, dimension grows to A1:Y36

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "covid19" PREFIX row at row 8 -------------------------
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "PREFIX"
$ws.Range("B8").Value = "covid19"
$ws.Range("C8").Value = "http://purl.org/zonmw/covid19/"

# --- Insert the new "spatial scope" concept row at row 19 -----------------
# (after the previous insert, the old row 18 "id-amr:10000/Deprecated" now
# sits at row 18 still, and old row 19's "id-amr:10001/diary cattle" sits at
# row 19 -- we insert above it so the new concept lands at row 19 and
# "diary cattle" moves to row 20)
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "covid19:10001"
$ws.Range("B19").Value = "spatial scope"

# --- "diary cattle" (now row 20) references the new concept as broader ----
$ws.Range("E20").Value = "covid19:10001"
$ws.Range("G20").Value = " "
